$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Table 1 (sheet1): add a new column I ("C") with saturation-pressure data,
# flip the sign on several G/H pairs (Reaction log10K / dH values).
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Table 1")

$ws1.Range("I1").Value = "C"
$ws1.Range("I2").Value = 2655395.9

$ws1.Range("G2").Value = 23.278899717000002
$ws1.Range("H2").Value = -99057.305063000007

$ws1.Range("G10").Value = -9.5500000000000007
$ws1.Range("H10").Value = 63948

$ws1.Range("G11").Value = -26.91
$ws1.Range("H11").Value = 204359

$ws1.Range("G12").Value = -29.86
$ws1.Range("H12").Value = 200903

$ws1.Range("H13").Value = -33554.101999999999

$ws1.Activate() | Out-Null
$ws1.Range("H14").Select() | Out-Null

# ---------------------------------------------------------------------------
# Table 2 (sheet2): add stoichiometric coefficient in F41.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Table 2")

$ws2.Range("F41").Value = 0.5

$ws2.Activate() | Out-Null
$ws2.Range("F41").Select() | Out-Null

# ---------------------------------------------------------------------------
# Table 3 (sheet3): refine the fitted coefficients for rows 2-3.
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Table 3")

$ws3.Range("C2").Value = -0.94840000000000002
$ws3.Range("D2").Value = 7404.5

$ws3.Range("C3").Value = 0.41520000000000001
$ws3.Range("D3").Value = 2330.6170000000002

$ws3.Activate() | Out-Null
$ws3.Range("C45").Select() | Out-Null

# ---------------------------------------------------------------------------
# Table 4 (sheet4): insert a new column E ("C" header / saturation data),
# pushing the old "Reactants" column to F; update row 3's C/D values and
# add a formula-driven saturation-pressure cell in E3.
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Table 4")

$ws4.Columns("E:E").Insert()

$ws4.Range("E1").Value = "C"

$ws4.Range("C3").Value = -12.8489
$ws4.Range("D3").Value = 43367.3
$ws4.Range("E3").Formula = "=-2655400"

$ws4.Activate() | Out-Null
$excel.ActiveWindow.Zoom = 130
$ws4.Range("F36").Select() | Out-Null
